$wb = $excel.ActiveWorkbook

# zh-cn sheet
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("D4").Value = "2016-01-19 04:16:33"
$wsZhCn.Range("G4").Value = "2016-01-19 04:17:28"

# de-de sheet
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("D4").Value = "2016-01-19 04:16:42"
$wsDeDe.Range("G4").Value = "2016-01-19 04:17:46"
